$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells from the crypto price/volume refresh.
# D (Price) and B/C (Coin/Link) columns are plain text in this sheet
# (values like "52.20" or "0.0800" must stay text, not become numbers),
# so force text format before writing them.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "42.949.97"
$ws.Cells.Item(2, 5).Value = "  -1.13%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.338.23"
$ws.Cells.Item(3, 5).Value = "  +1.12%  "

$ws.Cells.Item(4, 5).Value = "  +0.04%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "306.42"
$ws.Cells.Item(5, 5).Value = "  -1.68%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "101.51"
$ws.Cells.Item(6, 5).Value = "  -0.89%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.511"
$ws.Cells.Item(7, 5).Value = "  -4.67%  "

$ws.Cells.Item(8, 5).Value = "  +0.06%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.511"
$ws.Cells.Item(9, 5).Value = "  -3.27%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "35.13"
$ws.Cells.Item(10, 5).Value = "  -2.44%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "52.20"
$ws.Cells.Item(11, 5).Value = "  +1.16%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0800"
$ws.Cells.Item(12, 5).Value = "  -1.94%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.112"
$ws.Cells.Item(13, 5).Value = "  +0.04%  "

$ws.Cells.Item(14, 5).Value = "  -2.54%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.92"
$ws.Cells.Item(15, 5).Value = "  +5.89%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2.311.99"
$ws.Cells.Item(16, 5).Value = "  -0.22%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.810"
$ws.Cells.Item(17, 5).Value = "  -0.25%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "42.896.58"

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.21"
$ws.Cells.Item(19, 5).Value = "  +0.25%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0₃0911"
$ws.Cells.Item(20, 5).Value = "  -2.68%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.68"
$ws.Cells.Item(21, 5).Value = "  -5.64%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "67.95"
$ws.Cells.Item(22, 5).Value = "  -0.24%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "236.64"
$ws.Cells.Item(23, 5).Value = "  -2.07%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.02"
$ws.Cells.Item(24, 5).Value = "  +0.59%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.55"
$ws.Cells.Item(25, 5).Value = "  -2.78%  "

$ws.Cells.Item(26, 5).Value = "  -0.13%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "25.42"
$ws.Cells.Item(27, 5).Value = "  +3.04%  "

$ws.Cells.Item(28, 5).Value = "  +9.62%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "35.16"
$ws.Cells.Item(29, 5).Value = "  -5.23%  "

$ws.Cells.Item(30, 5).Value = "  -2.61%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "161.24"
$ws.Cells.Item(31, 5).Value = "  -4.19%  "

$ws.Cells.Item(32, 5).Value = "  +0.02%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.13"
$ws.Cells.Item(33, 5).Value = "  -3.00%  "

$ws.Cells.Item(34, 2).NumberFormat = "@"
$ws.Cells.Item(34, 3).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 2).Value = "WEMIXToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(34, 4).Value = "2.48"
$ws.Cells.Item(34, 5).Value = "  -1.00%  "

$ws.Cells.Item(35, 2).NumberFormat = "@"
$ws.Cells.Item(35, 3).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 2).Value = "RenderToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(35, 4).Value = "4.64"
$ws.Cells.Item(35, 5).Value = "  +6.41%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "17.54"
$ws.Cells.Item(36, 5).Value = "  -0.67%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0727"
$ws.Cells.Item(37, 5).Value = "  -2.28%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.93"
$ws.Cells.Item(38, 5).Value = "  -4.38%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.86"
$ws.Cells.Item(39, 5).Value = "  -1.67%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.102"
$ws.Cells.Item(40, 5).Value = "  -3.44%  "

$ws.Cells.Item(41, 5).Value = "  -2.34%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.47"
$ws.Cells.Item(42, 5).Value = "  +6.38%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.010.29"
$ws.Cells.Item(43, 5).Value = "  +1.82%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0285"
$ws.Cells.Item(44, 5).Value = "  -1.33%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "18.80"
$ws.Cells.Item(45, 5).Value = "  -3.17%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.26"
$ws.Cells.Item(46, 5).Value = "  +3.42%  "

$ws.Cells.Item(47, 5).Value = "  -1.42%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "55.83"
$ws.Cells.Item(48, 5).Value = "  +0.41%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.90"
$ws.Cells.Item(49, 5).Value = "  -0.56%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.565.76"
$ws.Cells.Item(50, 5).Value = "  +1.09%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "4.74"
$ws.Cells.Item(51, 5).Value = "  +3.23%  "
